$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "name" values (column B) to use underscore-separated naming
$ws.Range("B2").Value = "Good_Ending"
$ws.Range("B3").Value = "Bad_Ending"
$ws.Range("B4").Value = "Clear_Game_Once"

# Add new row 5 for the "All_Achievements" entry
$ws.Range("A5").Value = 99999
$ws.Range("B5").Value = "All_Achievements"
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = -1
$ws.Range("E5").Value = "allAchievement_99999"
$ws.Range("F5").Value = $false
$ws.Range("G5").Value = -1
$ws.Range("H5").Value = 0

# Update selection to match the new active cell
$ws.Range("H5").Select()
